# Apply the "as of" date roll + updated weight/percent-change figures to the
# COREINTL_holdings model-holdings sheet.
#
# The sheet carries legacy sheet protection (password "D382"), so cell writes
# must be bracketed by Unprotect/Protect.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect("D382")

# --- Update the "as of" date in the confidential disclosure note (A7) ---
$oldNote = $ws.Range("A7").Value2
$newNote = $oldNote -replace "2021-04-05", "2021-04-06"
$ws.Range("A7").Value = $newNote

# --- Update the weight / percent-change figures (rows 2-4) ---
$ws.Range("D2").Value = 0.8439515721727586
$ws.Range("E2").Value = -0.01064375480892543

$ws.Range("D3").Value = 0.1560484278272414
$ws.Range("E3").Value = 0.005548363232846354

$ws.Range("E4").Value = -0.00811700024531381

$ws.Protect("D382")
